$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Remove the first 4 data rows (spreadsheet rows 2-5); remaining data
    # rows shift up so that old row 6 becomes new row 2, etc.
    $ws.Range("A2:C5").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)

    # Renumber column A (the "Cutoff" index) 0..14 for the remaining 15
    # data rows (now in rows 2..16).
    for ($i = 0; $i -le 14; $i++) {
        $r = $i + 2
        $ws.Cells.Item($r, 1).Value = $i
    }
}
